$d = $word.ActiveDocument

$newText = "Campaign Dates that use Gemini: February 14-23, March 14-24"

# Collect the paragraphs that contain the outdated campaign-dates sentence.
# (Collecting indices first, since we mutate paragraph content as we go.)
$targets = New-Object System.Collections.ArrayList
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Campaign Dates that use Perseus*") {
        [void]$targets.Add($i)
    }
    $i = $i + 1
}

foreach ($idx in $targets) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    # Exclude the trailing paragraph mark from the range.
    $r.End = $r.End - 1
    # Clear every existing run (and its formatting) in the paragraph body.
    $r.Text = ""
    # Re-collapse to the (now empty) paragraph body and insert plain text,
    # which creates a single run with no explicit run formatting.
    $r2 = $p.Range
    $r2.End = $r2.End - 1
    $r2.InsertAfter($newText)
}
